$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links) - safe to assign directly,
# these values never look like numbers so Excel keeps them as text.
$textUpdates = @(
    @{ Cell = "B6"; Value = "FTXToken" },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "B7"; Value = "BTSEToken" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "B8"; Value = "MXToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "B9"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "B11"; Value = "MandalaExchangeToken" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "B12"; Value = "BitrueCoin" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "B13"; Value = "BitMartToken" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "B14"; Value = "BitForexToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "B15"; Value = "TigerCash" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "B16"; Value = "LEO" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "B17"; Value = "GateToken" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Numeric-looking text updates (price/volume columns) - these must stay
# stored as literal text (matching the source inlineStr cells), so force
# the cell to Text format before writing, then drop back to the default
# "Normal" style so no stray number-format style is left attached.
$numericTextUpdates = @(
    @{ Cell = "D2"; Value = "308.64" },
    @{ Cell = "E2"; Value = "0.51%" },
    @{ Cell = "D3"; Value = "40.90" },
    @{ Cell = "E3"; Value = "0.48%" },
    @{ Cell = "D4"; Value = "5.129" },
    @{ Cell = "E4"; Value = "1.92%" },
    @{ Cell = "D5"; Value = "0.07624" },
    @{ Cell = "E5"; Value = "0.23%" },
    @{ Cell = "D6"; Value = "1.619" },
    @{ Cell = "E6"; Value = "1.61%" },
    @{ Cell = "D7"; Value = "2.488" },
    @{ Cell = "E7"; Value = "2.84%" },
    @{ Cell = "D8"; Value = "0.9094" },
    @{ Cell = "E8"; Value = "0.59%" },
    @{ Cell = "D9"; Value = "0.1238" },
    @{ Cell = "E9"; Value = "26.59%" },
    @{ Cell = "D10"; Value = "0.1813" },
    @{ Cell = "E10"; Value = "2.81%" },
    @{ Cell = "D11"; Value = "0.09082" },
    @{ Cell = "E11"; Value = "-1.16%" },
    @{ Cell = "D12"; Value = "0.04276" },
    @{ Cell = "E12"; Value = "-1.07%" },
    @{ Cell = "D13"; Value = "0.1046" },
    @{ Cell = "E13"; Value = "-0.56%" },
    @{ Cell = "D14"; Value = "0.001260" },
    @{ Cell = "E14"; Value = "1.94%" },
    @{ Cell = "D15"; Value = "0.005861" },
    @{ Cell = "E15"; Value = "0.74%" },
    @{ Cell = "D16"; Value = "3.351" },
    @{ Cell = "E16"; Value = "-0.59%" },
    @{ Cell = "D17"; Value = "4.285" },
    @{ Cell = "E17"; Value = "1.03%" },
    @{ Cell = "E18"; Value = "-0.67%" },
    @{ Cell = "D19"; Value = "6.912" },
    @{ Cell = "E19"; Value = "1.21%" },
    @{ Cell = "D20"; Value = "0.1395" },
    @{ Cell = "E20"; Value = "3.24%" },
    @{ Cell = "E21"; Value = "-0.69%" },
    @{ Cell = "D22"; Value = "0.04049" },
    @{ Cell = "E22"; Value = "-2.65%" },
    @{ Cell = "E23"; Value = "4.61%" },
    @{ Cell = "D24"; Value = "0.004061" },
    @{ Cell = "E24"; Value = "-0.06%" },
    @{ Cell = "D25"; Value = "0.0001272" },
    @{ Cell = "E25"; Value = "-2.19%" },
    @{ Cell = "E26"; Value = "24.66%" },
    @{ Cell = "D38"; Value = "0.02421" },
    @{ Cell = "E38"; Value = "-0.08%" },
    @{ Cell = "E39"; Value = "1.67%" },
    @{ Cell = "D40"; Value = "0.007831" },
    @{ Cell = "E40"; Value = "-0.08%" },
    @{ Cell = "E41"; Value = "-0.20%" },
    @{ Cell = "D42"; Value = "0.006805" },
    @{ Cell = "E42"; Value = "-3.82%" },
    @{ Cell = "D43"; Value = "0.001933" },
    @{ Cell = "E43"; Value = "-0.82%" },
    @{ Cell = "D44"; Value = "0.008077" },
    @{ Cell = "E44"; Value = "-3.45%" },
    @{ Cell = "D45"; Value = "0.3065" },
    @{ Cell = "E45"; Value = "-7.78%" },
    @{ Cell = "D46"; Value = "0.00006898" },
    @{ Cell = "E46"; Value = "7.20%" },
    @{ Cell = "E47"; Value = "0.11%" },
    @{ Cell = "D48"; Value = "0.1072" },
    @{ Cell = "E48"; Value = "1,754.03%" },
    @{ Cell = "E50"; Value = "0.11%" },
    @{ Cell = "E51"; Value = "0.11%" }
)

foreach ($u in $numericTextUpdates) {
    $cellRange = $ws.Range($u.Cell)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $u.Value
    $cellRange.Style = "Normal"
}
